$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5885016
$ws.Range("I32").Value = 1250
$ws.Range("J32").Value = 6669518.5
$ws.Range("K32").Value = 1250
$ws.Range("L32").Value = 6669518.5
$ws.Range("M32").Value = -924
$ws.Range("N32").Value = -6670170.5
$ws.Range("H51").Value = 3696.8
$ws.Range("I51").Value = 3494.3333
$ws.Range("J51").Value = 4000.5
$ws.Range("K51").Value = 3494.3333
$ws.Range("L51").Value = 4000.5
$ws.Range("M51").Value = -3010.3333
$ws.Range("N51").Value = -4968.5
$ws.Range("H80").Value = 1080.6786
$ws.Range("I80").Value = 717.9167
$ws.Range("J80").Value = 1352.75
$ws.Range("K80").Value = 2153.7501
$ws.Range("L80").Value = 4058.25
$ws.Range("M80").Value = -1155.7501
$ws.Range("N80").Value = -6054.25
$ws.Range("H83").Value = 1080.6786
$ws.Range("I83").Value = 717.9167
$ws.Range("J83").Value = 1352.75
$ws.Range("K83").Value = 6461.2503
$ws.Range("L83").Value = 12174.75
$ws.Range("M83").Value = -1469.2503
$ws.Range("N83").Value = -22158.75
$ws.Range("H88").Value = 4897.6875
$ws.Range("I88").Value = 5645.125
$ws.Range("K88").Value = 5645.125
$ws.Range("M88").Value = -5239.125
$ws.Range("H91").Value = 4897.6875
$ws.Range("I91").Value = 5645.125
$ws.Range("K91").Value = 5645.125
$ws.Range("M91").Value = -4241.125
$ws.Range("H92").Value = 595.8570999999999
$ws.Range("J92").Value = 1111.75
$ws.Range("L92").Value = 1111.75
$ws.Range("N92").Value = -3607.75
$ws.Range("H94").Value = 3249.5
$ws.Range("I94").Value = 3249.5
$ws.Range("K94").Value = 3249.5
$ws.Range("M94").Value = -2798.5
$ws.Range("H98").Value = 1004.4167
$ws.Range("I98").Value = 1004.6087
$ws.Range("K98").Value = 1004.6087
$ws.Range("M98").Value = 493.3913
$ws.Range("H122").Value = 1004.4167
$ws.Range("I122").Value = 1004.6087
$ws.Range("K122").Value = 3013.8261
$ws.Range("M122").Value = -563.8261000000002
$ws.Range("H125").Value = 3907
$ws.Range("I125").Value = 1160.75
$ws.Range("K125").Value = 10446.75
$ws.Range("M125").Value = -7986.75
$ws.Range("H128").Value = 78076.92
$ws.Range("J128").Value = 78076.92
$ws.Range("L128").Value = 78076.92
$ws.Range("N128").Value = -88036.92
$ws.Range("H131").Value = 6797
$ws.Range("I131").Value = 4121.25
$ws.Range("K131").Value = 12363.75
$ws.Range("M131").Value = -7323.75
$ws.Range("H132").Value = 23259.107
$ws.Range("I132").Value = 3862.45
$ws.Range("K132").Value = 11587.35
$ws.Range("M132").Value = -9057.349999999999
$ws.Range("H137").Value = 2505
$ws.Range("I137").Value = 2396.4666
$ws.Range("K137").Value = 7189.399800000001
$ws.Range("M137").Value = -4639.399800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12992843
$ws.Range("I32").Value = 13519336
$ws.Range("K32").Value = 13519336
$ws.Range("M32").Value = -13519049
$ws.Range("H45").Value = 6008
$ws.Range("I45").Value = 4070
$ws.Range("J45").Value = 9238
$ws.Range("K45").Value = 4070
$ws.Range("L45").Value = 9238
$ws.Range("M45").Value = -3693
$ws.Range("N45").Value = -9992
$ws.Range("H97").Value = 1276.4375
$ws.Range("I97").Value = 672.4516
$ws.Range("K97").Value = 672.4516
$ws.Range("M97").Value = -176.4516
$ws.Range("H102").Value = 2023.6
$ws.Range("I102").Value = 2023.6
$ws.Range("K102").Value = 2023.6
$ws.Range("M102").Value = -401.5999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2374.4688
$ws.Range("I20").Value = 1160.6111
$ws.Range("J20").Value = 3935.1428
$ws.Range("K20").Value = 1160.6111
$ws.Range("L20").Value = 3935.1428
$ws.Range("M20").Value = -913.6111000000001
$ws.Range("N20").Value = -4429.1428
$ws.Range("H99").Value = 1470.9286
$ws.Range("I99").Value = 1363
$ws.Range("K99").Value = 1363
$ws.Range("M99").Value = 135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11711
$ws.Range("I16").Value = 20731.666
$ws.Range("J16").Value = 6298.6
$ws.Range("K16").Value = 20731.666
$ws.Range("L16").Value = 6298.6
$ws.Range("M16").Value = -20444.666
$ws.Range("N16").Value = -6872.6
$ws.Range("H31").Value = 2081
$ws.Range("I31").Value = 1997.6428
$ws.Range("J31").Value = 2340.3333
$ws.Range("K31").Value = 1997.6428
$ws.Range("L31").Value = 2340.3333
$ws.Range("M31").Value = -1702.6428
$ws.Range("N31").Value = -2930.3333
$ws.Range("H34").Value = 2081
$ws.Range("I34").Value = 1997.6428
$ws.Range("J34").Value = 2340.3333
$ws.Range("K34").Value = 1997.6428
$ws.Range("L34").Value = 2340.3333
$ws.Range("M34").Value = -1795.6428
$ws.Range("N34").Value = -2744.3333
$ws.Range("H58").Value = 2829.3076
$ws.Range("J58").Value = 5750
$ws.Range("L58").Value = 5750
$ws.Range("N58").Value = -6156
$ws.Range("H99").Value = 14965121
$ws.Range("I99").Value = 4880788
$ws.Range("K99").Value = 4880788
$ws.Range("M99").Value = -4879290
$ws.Range("H113").Value = 11711
$ws.Range("I113").Value = 20731.666
$ws.Range("J113").Value = 6298.6
$ws.Range("K113").Value = 20731.666
$ws.Range("L113").Value = 6298.6
$ws.Range("M113").Value = -18561.666
$ws.Range("N113").Value = -10638.6
$ws.Range("H126").Value = 14965121
$ws.Range("I126").Value = 4880788
$ws.Range("K126").Value = 14642364
$ws.Range("M126").Value = -14639894
$ws.Range("H132").Value = 2420.8667
$ws.Range("I132").Value = 2450.9285
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7352.7855
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -4822.7855
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 2139.6858
$ws.Range("I134").Value = 2153.4243
$ws.Range("J134").Value = 1913
$ws.Range("K134").Value = 6460.2729
$ws.Range("L134").Value = 5739
$ws.Range("M134").Value = -3925.2729
$ws.Range("N134").Value = -10809
$ws.Range("H136").Value = 2829.3076
$ws.Range("J136").Value = 5750
$ws.Range("L136").Value = 17250
$ws.Range("N136").Value = -22350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 2009.8
$ws.Range("I8").Value = 2009.8
$ws.Range("K8").Value = 6029.4
$ws.Range("M8").Value = -5890.4
$ws.Range("H59").Value = 108666.664
$ws.Range("J59").Value = 1000
$ws.Range("L59").Value = 3000
$ws.Range("N59").Value = -4080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1578822.1
$ws.Range("I3").Value = 250375
$ws.Range("J3").Value = 3350085
$ws.Range("K3").Value = 250375
$ws.Range("L3").Value = 3350085
$ws.Range("M3").Value = -250259
$ws.Range("N3").Value = -3350317
$ws.Range("H10").Value = 647.6
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = 797
$ws.Range("K10").Value = 50
$ws.Range("L10").Value = 797
$ws.Range("M10").Value = 119
$ws.Range("N10").Value = -1135
$ws.Range("H97").Value = 690.7059
$ws.Range("I97").Value = 729.1
$ws.Range("J97").Value = 635.8570999999999
$ws.Range("K97").Value = 729.1
$ws.Range("L97").Value = 635.8570999999999
$ws.Range("M97").Value = -233.1
$ws.Range("N97").Value = -1627.8571
$ws.Range("H99").Value = 13449.333
$ws.Range("I99").Value = 13449.333
$ws.Range("K99").Value = 13449.333
$ws.Range("M99").Value = -11203.333
$ws.Range("H102").Value = 4416.727
$ws.Range("I102").Value = 3170.375
$ws.Range("K102").Value = 3170.375
$ws.Range("M102").Value = -1548.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8159.2354
$ws.Range("I7").Value = 8686.546
$ws.Range("K7").Value = 8686.546
$ws.Range("M7").Value = -8574.546
$ws.Range("H22").Value = 2200
$ws.Range("J22").Value = 2316.6667
$ws.Range("L22").Value = 2316.6667
$ws.Range("N22").Value = -2906.6667
$ws.Range("H27").Value = 2200
$ws.Range("J27").Value = 2316.6667
$ws.Range("L27").Value = 2316.6667
$ws.Range("N27").Value = -2530.6667
$ws.Range("H46").Value = 1483.1666
$ws.Range("I46").Value = 1000.6667
$ws.Range("J46").Value = 1965.6666
$ws.Range("K46").Value = 1000.6667
$ws.Range("L46").Value = 1965.6666
$ws.Range("M46").Value = -812.6667
$ws.Range("N46").Value = -2341.6666
$ws.Range("H55").Value = 3186.3125
$ws.Range("I55").Value = 2983.75
$ws.Range("J55").Value = 3388.875
$ws.Range("K55").Value = 2983.75
$ws.Range("L55").Value = 3388.875
$ws.Range("M55").Value = -2810.75
$ws.Range("N55").Value = -3734.875
$ws.Range("H126").Value = 8159.2354
$ws.Range("I126").Value = 8686.546
$ws.Range("K126").Value = 26059.638
$ws.Range("M126").Value = -23589.638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1254.1428
$ws.Range("I113").Value = 698.44446
$ws.Range("J113").Value = 2254.4
$ws.Range("K113").Value = 2095.33338
$ws.Range("L113").Value = 6763.200000000001
$ws.Range("M113").Value = 74.66661999999997
$ws.Range("N113").Value = -11103.2
$ws.Range("H132").Value = 1985.5555
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 69950
$ws.Range("L137").Value = 69950
$ws.Range("N137").Value = -80150
